# ============================================================
# Sheet "Results": add TBT, CLS, GreenHostName columns; refresh
# CO2 calc values (median-based Lighthouse aggregation).
# ============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Insert new columns: G:H (TBT, CLS) before PageWeightBytes,
# and a new GreenHostName column after GreenHosting (old col L).
$ws.Range("G1:H1").EntireColumn.Insert()
$ws.Range("M1").EntireColumn.Insert()

# New column headers
$ws.Cells.Item(1,7).Value = "TBT"
$ws.Cells.Item(1,8).Value = "CLS"
$ws.Cells.Item(1,13).Value = "GreenHostName"

# Row 2:
$ws.Cells.Item(2,7).Value = 252.0
$ws.Cells.Item(2,8).Value = 0.0
$ws.Cells.Item(2,12).Value = "No"
$ws.Cells.Item(2,15).Value = 2.170902756648601
$ws.Cells.Item(2,16).Value = 1.648022439537

# Row 3:
$ws.Cells.Item(3,7).Value = 5.5
$ws.Cells.Item(3,8).Value = 0.0
$ws.Cells.Item(3,12).Value = "No"
$ws.Cells.Item(3,15).Value = 0.5888359522656
$ws.Cells.Item(3,16).Value = 0.4470098255519999

# Row 4:
$ws.Cells.Item(4,7).Value = 58.0
$ws.Cells.Item(4,8).Value = 0.0068852206501101
$ws.Cells.Item(4,12).Value = "Yes"
$ws.Cells.Item(4,13).Value = "Cloudflare"
$ws.Cells.Item(4,15).Value = 0.4535404126578
$ws.Cells.Item(4,16).Value = 0.344301362651

# Row 5:
$ws.Cells.Item(5,7).Value = 77.0
$ws.Cells.Item(5,8).Value = 0.0
$ws.Cells.Item(5,12).Value = "Yes"
$ws.Cells.Item(5,13).Value = "Cloudflare"
$ws.Cells.Item(5,15).Value = 0.4541828401656
$ws.Cells.Item(5,16).Value = 0.344789056052

# Row 6:
$ws.Cells.Item(6,7).Value = 177.5
$ws.Cells.Item(6,8).Value = 0.0056225658040376
$ws.Cells.Item(6,12).Value = "No"
$ws.Cells.Item(6,15).Value = 0.548921168874
$ws.Cells.Item(6,16).Value = 0.4167088558299999

# Row 7:
$ws.Cells.Item(7,7).Value = 71.5
$ws.Cells.Item(7,8).Value = 0.0
$ws.Cells.Item(7,12).Value = "No"
$ws.Cells.Item(7,15).Value = 1.2336768728856
$ws.Cells.Item(7,16).Value = 0.936535348452

# Row 8:
$ws.Cells.Item(8,7).Value = 144.0
$ws.Cells.Item(8,8).Value = 0.0
$ws.Cells.Item(8,12).Value = "No"
$ws.Cells.Item(8,15).Value = 0.4257275537736
$ws.Cells.Item(8,16).Value = 0.3231874664119999

# Row 9:
$ws.Cells.Item(9,7).Value = 194.5000000000009
$ws.Cells.Item(9,8).Value = 0.0
$ws.Cells.Item(9,12).Value = "No"
$ws.Cells.Item(9,15).Value = 0.5226678343692001
$ws.Cells.Item(9,16).Value = 0.3967788593139999

# Row 10:
$ws.Cells.Item(10,7).Value = 175.5
$ws.Cells.Item(10,8).Value = 0.0
$ws.Cells.Item(10,12).Value = "No"
$ws.Cells.Item(10,15).Value = 0.5810038530426
$ws.Cells.Item(10,16).Value = 0.4410641537669999

# Row 11:
$ws.Cells.Item(11,7).Value = 30.5
$ws.Cells.Item(11,8).Value = 0.0007589648594464
$ws.Cells.Item(11,12).Value = "Yes"
$ws.Cells.Item(11,13).Value = "Sitevision AB"
$ws.Cells.Item(11,15).Value = 1.192432797036
$ws.Cells.Item(11,16).Value = 0.90522525762

# Row 12:
$ws.Cells.Item(12,7).Value = 3.5
$ws.Cells.Item(12,8).Value = 0.0592427032813799
$ws.Cells.Item(12,12).Value = "No"
$ws.Cells.Item(12,15).Value = 0.3556255729806
$ws.Cells.Item(12,16).Value = 0.269970141477

# Row 13:
$ws.Cells.Item(13,7).Value = 131.0
$ws.Cells.Item(13,8).Value = 0.0
$ws.Cells.Item(13,12).Value = "No"
$ws.Cells.Item(13,15).Value = 0.9792456994476
$ws.Cells.Item(13,16).Value = 0.7433860782419999

# Row 14:
$ws.Cells.Item(14,7).Value = 2.0
$ws.Cells.Item(14,8).Value = 0.0038646982166687
$ws.Cells.Item(14,12).Value = "Yes"
$ws.Cells.Item(14,13).Value = "Sitevision AB"
$ws.Cells.Item(14,15).Value = 0.4746964660542
$ws.Cells.Item(14,16).Value = 0.3603618013889999

# Row 15:
$ws.Cells.Item(15,7).Value = 127.5000000000009
$ws.Cells.Item(15,8).Value = 0.0
$ws.Cells.Item(15,12).Value = "No"
$ws.Cells.Item(15,15).Value = 0.5742022427856
$ws.Cells.Item(15,16).Value = 0.4359007689519999

# Row 16:
$ws.Cells.Item(16,7).Value = 120.5
$ws.Cells.Item(16,8).Value = 0.0
$ws.Cells.Item(16,12).Value = "No"
$ws.Cells.Item(16,15).Value = 0.6544922734116001
$ws.Cells.Item(16,16).Value = 0.4968522656219999

# Row 17:
$ws.Cells.Item(17,7).Value = 0.0
$ws.Cells.Item(17,8).Value = 0.0
$ws.Cells.Item(17,12).Value = "No"
$ws.Cells.Item(17,15).Value = 2.3897701852362
$ws.Cells.Item(17,16).Value = 1.814173794079

# Row 18:
$ws.Cells.Item(18,7).Value = 133.0
$ws.Cells.Item(18,8).Value = 0.0123655455877825
$ws.Cells.Item(18,12).Value = "No"
$ws.Cells.Item(18,15).Value = 0.9521767846422002
$ws.Cells.Item(18,16).Value = 0.7228369408489999

# Row 19:
$ws.Cells.Item(19,7).Value = 287.5
$ws.Cells.Item(19,8).Value = 0.00003829323822311932
$ws.Cells.Item(19,12).Value = "No"
$ws.Cells.Item(19,15).Value = 1.210089018762
$ws.Cells.Item(19,16).Value = 0.91862882879

# Row 20:
$ws.Cells.Item(20,7).Value = 90.49999999999818
$ws.Cells.Item(20,8).Value = 0.0
$ws.Cells.Item(20,12).Value = "No"
$ws.Cells.Item(20,15).Value = 2.586074502445201
$ws.Cells.Item(20,16).Value = 1.963196553734

# Row 21:
$ws.Cells.Item(21,7).Value = 61.5
$ws.Cells.Item(21,8).Value = 0.0
$ws.Cells.Item(21,12).Value = "No"
$ws.Cells.Item(21,15).Value = 0.6110209623024
$ws.Cells.Item(21,16).Value = 0.4638513880079999

# Row 22:
$ws.Cells.Item(22,7).Value = 49.99999999999909
$ws.Cells.Item(22,8).Value = 0.0
$ws.Cells.Item(22,12).Value = "Yes"
$ws.Cells.Item(22,13).Value = "Sitevision AB"
$ws.Cells.Item(22,15).Value = 0.332777065959
$ws.Cells.Item(22,16).Value = 0.2526248909049999

# ============================================================
# Sheet "Averages": insert Avg_TBT_ms and Avg_CLS rows after
# Avg_SpeedIndex_ms (old row 5); refresh CO2 averages.
# ============================================================
$wsAvg = $wb.Worksheets.Item("Averages")

# Insert two new rows right before the old "Avg_Requests" row (row 6)
$wsAvg.Range("A6:B7").EntireRow.Insert()

$wsAvg.Cells.Item(6,1).Value = "Avg_TBT_ms"
$wsAvg.Cells.Item(6,2).Value = 104.43
$wsAvg.Cells.Item(7,1).Value = "Avg_CLS"
$wsAvg.Cells.Item(7,2).Value = 0.004

# Refresh the CO2 averages to reflect the recalculated per-site values
$wsAvg.Cells.Item(11,2).Value = 0.9187
$wsAvg.Cells.Item(12,2).Value = 0.6974
